$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 206; existing rows 206-246 shift down to 207-247.
$ws.Rows("206").Insert()

# Populate the newly inserted row 206 with the new data record.
$ws.Range("A206").Value = 7
$ws.Range("B206").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C206").Value = "Ñuble"
$ws.Range("D206").Value = 44995
$ws.Range("E206").Value = 16
$ws.Range("F206").Value = 100112045
$ws.Range("G206").Value = "Zapallo"
$ws.Range("H206").Value = "Camote"
$ws.Range("I206").Value = "1a (cosecha)"
$ws.Range("J206").Value = 200
$ws.Range("K206").Value = 350
$ws.Range("L206").Value = 400
$ws.Range("M206").Value = 375
$ws.Range("N206").Value = "$/kilo (volumen en unidades)"
$ws.Range("O206").Value = "Región del Maule"
$ws.Range("P206").Value = 375
$ws.Range("Q206").Value = 1
$ws.Range("R206").Value = "Hortaliza"
